# Add data for 2022-09-07
# Updates year-to-date violent crime counts across the "Citywide Totals",
# "By Neighborhood" and several per-neighborhood sheets to reflect a newly
# ingested day of crime records (and a couple of late-arriving records for
# prior years).

$wb = $excel.ActiveWorkbook

function Set-Cells {
    param(
        [string]$SheetName,
        [hashtable]$Values
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($cell in $Values.Keys) {
        $ws.Range($cell).Value = $Values[$cell]
    }
}

Set-Cells -SheetName "Citywide Totals" -Values @{
    "C2" = 49
    "I2" = 83
    "B3" = 61
    "F3" = 94
    "G3" = 91
    "H3" = 94
    "I3" = 142
    "D6" = 306
    "E6" = 310
    "F6" = 391
    "G6" = 352
    "H6" = 319
    "I6" = 377
    "B7" = 372
    "C7" = 462
    "D7" = 474
    "E7" = 467
    "F7" = 553
    "G7" = 511
    "H7" = 503
    "I7" = 622
}

Set-Cells -SheetName "By Neighborhood" -Values @{
    "G8"  = 24
    "I21" = 14
    "E22" = 2
    "G28" = 32
    "B36" = 20
    "F36" = 37
    "H36" = 30
    "I36" = 30
    "H43" = 3
    "I47" = 21
    "D50" = 10
    "I61" = 3
    "E70" = 10
    "G74" = 13
    "G76" = 11
    "H77" = 21
    "C86" = 10
    "F86" = 8
    "B98" = 372
    "C98" = 462
    "D98" = 474
    "E98" = 467
    "F98" = 553
    "G98" = 511
    "H98" = 503
    "I98" = 622
}

Set-Cells -SheetName "Rogers Park" -Values @{
    "G5" = 6
    "G6" = 11
}

Set-Cells -SheetName "Roseland" -Values @{
    "H3" = 5
    "H7" = 21
}

Set-Cells -SheetName "Austin" -Values @{
    "G3" = 5
    "G6" = 24
}

Set-Cells -SheetName "Chinatown" -Values @{
    "I3" = 5
    "I7" = 14
}

Set-Cells -SheetName "Grand Crossing" -Values @{
    "B3" = 3
    "F3" = 15
    "H6" = 18
    "I6" = 14
    "B7" = 20
    "F7" = 37
    "H7" = 30
    "I7" = 30
}

Set-Cells -SheetName "Little Italy, UIC" -Values @{
    "D5" = 8
    "D6" = 10
}

Set-Cells -SheetName "Uptown" -Values @{
    "C2" = 3
    "F5" = 5
    "C6" = 10
    "F6" = 8
}

Set-Cells -SheetName "Englewood" -Values @{
    "G3" = 4
    "G7" = 32
}

Set-Cells -SheetName "River North" -Values @{
    "G5" = 11
    "G6" = 13
}

Set-Cells -SheetName "Irving Park" -Values @{
    "H5" = 3
    "H6" = 3
}

Set-Cells -SheetName "Lake View" -Values @{
    "I2" = 2
    "I6" = 21
}

Set-Cells -SheetName "Old Town" -Values @{
    "E4" = 9
    "E5" = 10
}

Set-Cells -SheetName "Clearing" -Values @{
    "C4" = 1
    "C5" = 2
}

Write-Host "Applied 2022-09-07 data update"
